$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update faturamento (revenue) values in column B for rows 2-63
# as per the updated ADD faturamento diario data

$ws.Cells.Item(2, 2).Value = 36711.94
$ws.Cells.Item(3, 2).Value = 48583.35
$ws.Cells.Item(4, 2).Value = 23366.39
$ws.Cells.Item(5, 2).Value = 34512.78
$ws.Cells.Item(6, 2).Value = 55688.46
$ws.Cells.Item(7, 2).Value = 39651.05
$ws.Cells.Item(8, 2).Value = 7727.37
$ws.Cells.Item(9, 2).Value = 2416.2
$ws.Cells.Item(10, 2).Value = 7630.52
$ws.Cells.Item(11, 2).Value = 19121.19
$ws.Cells.Item(12, 2).Value = 18067.5
$ws.Cells.Item(13, 2).Value = 14391.72
$ws.Cells.Item(14, 2).Value = 62436.85
$ws.Cells.Item(15, 2).Value = 19623.66
$ws.Cells.Item(16, 2).Value = 22976.52
$ws.Cells.Item(17, 2).Value = 21772.35
$ws.Cells.Item(18, 2).Value = 28370.31
$ws.Cells.Item(19, 2).Value = 25609.93
$ws.Cells.Item(20, 2).Value = 19005
$ws.Cells.Item(21, 2).Value = 8316.389999999999
$ws.Cells.Item(22, 2).Value = 56281.38
$ws.Cells.Item(23, 2).Value = 62723.95
$ws.Cells.Item(24, 2).Value = 50954.74
$ws.Cells.Item(25, 2).Value = 21522.57
$ws.Cells.Item(26, 2).Value = 38519.31
$ws.Cells.Item(27, 2).Value = 13587.25
$ws.Cells.Item(28, 2).Value = 21015.05
$ws.Cells.Item(29, 2).Value = 10544.77
$ws.Cells.Item(30, 2).Value = 29469.84
$ws.Cells.Item(31, 2).Value = 14291
$ws.Cells.Item(32, 2).Value = 21396.45
$ws.Cells.Item(33, 2).Value = 26250.86
$ws.Cells.Item(34, 2).Value = 4035.97
$ws.Cells.Item(35, 2).Value = 19913.85
$ws.Cells.Item(36, 2).Value = 5325.7
$ws.Cells.Item(37, 2).Value = 22282.62
$ws.Cells.Item(38, 2).Value = 33108.32
$ws.Cells.Item(39, 2).Value = 29400.36
$ws.Cells.Item(40, 2).Value = 19290.09
$ws.Cells.Item(41, 2).Value = 13162.5
$ws.Cells.Item(42, 2).Value = 11137.4
$ws.Cells.Item(43, 2).Value = 42434.74
$ws.Cells.Item(44, 2).Value = 29976.02
$ws.Cells.Item(45, 2).Value = 20469.43
$ws.Cells.Item(46, 2).Value = 21525.35
$ws.Cells.Item(47, 2).Value = 29665.53
$ws.Cells.Item(48, 2).Value = 22654.48
$ws.Cells.Item(49, 2).Value = 27794.15
$ws.Cells.Item(50, 2).Value = 8942.360000000001
$ws.Cells.Item(51, 2).Value = 26486.27
$ws.Cells.Item(52, 2).Value = 9704.469999999999
$ws.Cells.Item(53, 2).Value = 5815.6
$ws.Cells.Item(54, 2).Value = 8184.83
$ws.Cells.Item(55, 2).Value = 18780.45
$ws.Cells.Item(56, 2).Value = 23732.7
$ws.Cells.Item(57, 2).Value = 7291.1
$ws.Cells.Item(58, 2).Value = 9459.75
$ws.Cells.Item(59, 2).Value = 21801.33
$ws.Cells.Item(60, 2).Value = 12669.84
$ws.Cells.Item(61, 2).Value = 57238.78
$ws.Cells.Item(62, 2).Value = 42200.07
$ws.Cells.Item(63, 2).Value = 95702.25999999999
